# Apply the edits described by the diff:
#  - Add a new "Rho_Intercept Table" sheet (sheetId 3) at the end of the workbook
#    containing the Rho_Intercept naming/config lookup table.
#  - Re-style the "Species" label + header rows on the "CVs" sheet (bold label,
#    the percent value cells centred instead of right-aligned).
#  - Update the active-sheet/selection bookkeeping to reflect the new sheet
#    being the one in front, and CVs' own selection moving off of B18:H22.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. CVs sheet: move the "Species" label from A4/A17 up into A3/A16 (bold,
#    left-aligned) and shift the two small header rows down to match, leaving
#    B4:H4 / B17:H17 with only the numeric knot counts. Also re-centre the
#    percent data cells (B5:H13 / B18:H22) instead of right-aligning them.
# ---------------------------------------------------------------------------
$cvs = $wb.Worksheets.Item("CVs")

# --- Table 1 (rows 2-13) ---
$cvs.Range("A4").Value = $cvs.Range("A4").Value   # no-op placeholder not used
$speciesLabel1 = $cvs.Range("A4").Value
$cvs.Range("A4").ClearContents()
$cvs.Range("A3").Value = $speciesLabel1
$cvs.Range("A3").Font.Bold = $true

# --- Table 2 (rows 15-22) ---
$speciesLabel2 = $cvs.Range("A17").Value
$cvs.Range("A17").ClearContents()
$cvs.Range("A16").Value = $speciesLabel2
$cvs.Range("A16").Font.Bold = $true

# Re-centre (instead of right-align) the percent-formatted data cells.
$cvs.Range("B5:H13").HorizontalAlignment = -4108   # xlCenter
$cvs.Range("B18:H22").HorizontalAlignment = -4108  # xlCenter

# CVs' own selection is no longer on the old B18:H22 block.
$cvs.Range("J8").Select()

# ---------------------------------------------------------------------------
# 2. Add the new "Rho_Intercept Table" sheet at the end of the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$rho = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$rho.Name = "Rho_Intercept Table"

# Column A: Rho_Intercept short names.
$colA = @("Rho_Intercept Name", "FE", "RW-FE", "FE-RW", "RW", "AR-FE", "FE-AR", "AR")
for ($r = 0; $r -lt $colA.Length; $r++) {
    $rho.Cells.Item($r + 1, 1).Value = $colA[$r]
}

# Header row for columns B & C.
$rho.Cells.Item(1, 2).Value = "Encounter Probability"
$rho.Cells.Item(1, 3).Value = "Positive Catch Rate"

# Column B body (Encounter Probability component).
$colB = @("Fixed effect", "Random walk", "Fixed effect", "Random walk", "Autoregressive (lag-1)", "Fixed effect", "Autoregressive (lag-1)")
for ($r = 0; $r -lt $colB.Length; $r++) {
    $rho.Cells.Item($r + 2, 2).Value = $colB[$r]
}

# Column C body (Positive Catch Rate component).
$colC = @("Fixed effect", "Fixed effect", "Random walk", "Random walk", "Fixed effect", "Autoregressive (lag-1)", "Autoregressive (lag-1)")
for ($r = 0; $r -lt $colC.Length; $r++) {
    $rho.Cells.Item($r + 2, 3).Value = $colC[$r]
}

# Bold header row, then best-fit the three columns to their content.
$rho.Range("A1:C1").Font.Bold = $true
$rho.Columns.Item(1).AutoFit() | Out-Null
$rho.Columns.Item(2).AutoFit() | Out-Null
$rho.Columns.Item(3).AutoFit() | Out-Null

# New sheet becomes the active/front tab with its own selection.
$rho.Select()
$rho.Range("G15").Select()
